$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.364.05"
$ws.Range("E2").Value = "  -2.13%  "
$ws.Range("D3").Value = "2.488.43"
$ws.Range("E3").Value = "  -2.35%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "565.58"
$ws.Range("E5").Value = "  -2.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.61"
$ws.Range("E6").Value = "  -3.57%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.510"
$ws.Range("E8").Value = "  -2.33%  "
$ws.Range("D9").Value = "2.486.99"
$ws.Range("E9").Value = "  -2.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.157"
$ws.Range("E10").Value = "  -4.33%  "
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.352"
$ws.Range("E12").Value = "  -1.93%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.89"
$ws.Range("E13").Value = "  -0.90%  "
$ws.Range("D14").Value = "2.937.34"
$ws.Range("E14").Value = "  -2.65%  "
$ws.Range("D15").Value = "69.196.96"
$ws.Range("E15").Value = "  -2.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000174"
$ws.Range("E16").Value = "  -3.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.26"
$ws.Range("E17").Value = "  -4.34%  "
$ws.Range("D18").Value = "2.493.03"
$ws.Range("E18").Value = "  -2.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.17"
$ws.Range("E19").Value = "  -3.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.35"
$ws.Range("E20").Value = "  -7.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "343.81"
$ws.Range("E21").Value = "  -2.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.85"
$ws.Range("E22").Value = "  -2.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.91"
$ws.Range("E23").Value = "  -5.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.35"
$ws.Range("E25").Value = "  -1.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.87"
$ws.Range("E26").Value = "  -4.38%  "
$ws.Range("D27").Value = "2.611.64"
$ws.Range("E27").Value = "  -3.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.64"
$ws.Range("E28").Value = "  -3.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").Value = "0.0₃0869"
$ws.Range("E30").Value = "  -5.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.64"
$ws.Range("E31").Value = "  -3.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "441.79"
$ws.Range("E32").Value = "  -5.79%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.18"
$ws.Range("E33").Value = "  -7.06%  "
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.70"
$ws.Range("E35").Value = "  -3.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "156.18"
$ws.Range("E36").Value = "  -0.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.113"
$ws.Range("E37").Value = "  -4.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.03"
$ws.Range("E38").Value = "  -0.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.08"
$ws.Range("E39").Value = "  -4.26%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.314"
$ws.Range("E41").Value = "  -2.37%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.58"
$ws.Range("E42").Value = "  -2.83%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.56"
$ws.Range("E43").Value = "  -5.95%  "
$ws.Range("E44").Value = "  -8.58%  "
$ws.Range("E45").Value = "  -8.71%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "138.26"
$ws.Range("E46").Value = "  -4.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.42"
$ws.Range("E47").Value = "  -3.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.511"
$ws.Range("E48").Value = "  -4.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0727"
$ws.Range("E49").Value = "  -1.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.570"
$ws.Range("E50").Value = "  -2.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0920"
$ws.Range("E51").Value = "  -1.47%  "
